$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3654.4546
$ws.Range("I51").Value = 3459.8
$ws.Range("J51").Value = 3816.6667
$ws.Range("K51").Value = 3459.8
$ws.Range("L51").Value = 3816.6667
$ws.Range("M51").Value = -2975.8
$ws.Range("N51").Value = -4784.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2598700.8
$ws.Range("I132").Value = 3040718
$ws.Range("K132").Value = 9122154
$ws.Range("M132").Value = -9119624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 536.5192
$ws.Range("I2").Value = 351.86206
$ws.Range("J2").Value = 769.34784
$ws.Range("K2").Value = 351.86206
$ws.Range("L2").Value = 769.34784
$ws.Range("M2").Value = -238.86206
$ws.Range("N2").Value = -995.34784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1462.875
$ws.Range("I45").Value = 1339.6923
$ws.Range("J45").Value = 1996.6666
$ws.Range("K45").Value = 1339.6923
$ws.Range("L45").Value = 1996.6666
$ws.Range("M45").Value = -962.6922999999999
$ws.Range("N45").Value = -2750.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3335000
$ws.Range("I63").Value = 10000000
$ws.Range("K63").Value = 10000000
$ws.Range("M63").Value = -9999314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3335000
$ws.Range("I66").Value = 10000000
$ws.Range("K66").Value = 50000000
$ws.Range("M66").Value = -49996568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1840.25
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 2407.3333
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 2407.3333
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -3399.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1371.7142
$ws.Range("I110").Value = 1449.25
$ws.Range("J110").Value = 906.5
$ws.Range("K110").Value = 1449.25
$ws.Range("L110").Value = 906.5
$ws.Range("M110").Value = 595.75
$ws.Range("N110").Value = -4996.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 536.5192
$ws.Range("I116").Value = 351.86206
$ws.Range("J116").Value = 769.34784
$ws.Range("K116").Value = 351.86206
$ws.Range("L116").Value = 769.34784
$ws.Range("M116").Value = 1942.13794
$ws.Range("N116").Value = -5357.34784

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 536.5192
$ws.Range("I3").Value = 351.86206
$ws.Range("J3").Value = 769.34784
$ws.Range("K3").Value = 351.86206
$ws.Range("L3").Value = 769.34784
$ws.Range("M3").Value = -237.86206
$ws.Range("N3").Value = -997.34784

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2245.4546
$ws.Range("I86").Value = 2414.2856
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 2414.2856
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -1291.2856
$ws.Range("N86").Value = -4196

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2245.4546
$ws.Range("I89").Value = 2414.2856
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 12071.428
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = -6455.428
$ws.Range("N89").Value = -20982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3500
$ws.Range("J99").Value = 2750
$ws.Range("K99").Value = 3500
$ws.Range("L99").Value = 2750
$ws.Range("M99").Value = -2002
$ws.Range("N99").Value = -5746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2491
$ws.Range("I107").Value = 2558.5715
$ws.Range("J107").Value = 2333.3333
$ws.Range("K107").Value = 2558.5715
$ws.Range("L107").Value = 2333.3333
$ws.Range("M107").Value = -638.5715
$ws.Range("N107").Value = -6173.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3476078
$ws.Range("I31").Value = 2606.2593
$ws.Range("J31").Value = 7941970
$ws.Range("K31").Value = 2606.2593
$ws.Range("L31").Value = 7941970
$ws.Range("M31").Value = -2311.2593
$ws.Range("N31").Value = -7942560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3476078
$ws.Range("I34").Value = 2606.2593
$ws.Range("J34").Value = 7941970
$ws.Range("K34").Value = 2606.2593
$ws.Range("L34").Value = 7941970
$ws.Range("M34").Value = -2404.2593
$ws.Range("N34").Value = -7942374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1238.5834
$ws.Range("I58").Value = 1069.5
$ws.Range("J58").Value = 1576.75
$ws.Range("K58").Value = 1069.5
$ws.Range("L58").Value = 1576.75
$ws.Range("M58").Value = -866.5
$ws.Range("N58").Value = -1982.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 674.2222
$ws.Range("I105").Value = 497.5
$ws.Range("J105").Value = 895.125
$ws.Range("K105").Value = 497.5
$ws.Range("L105").Value = 895.125
$ws.Range("M105").Value = 1249.5
$ws.Range("N105").Value = -4389.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1238.5834
$ws.Range("I136").Value = 1069.5
$ws.Range("J136").Value = 1576.75
$ws.Range("K136").Value = 3208.5
$ws.Range("L136").Value = 4730.25
$ws.Range("M136").Value = -658.5
$ws.Range("N136").Value = -9830.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 200
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 738.26
$ws.Range("I131").Value = 301.81818
$ws.Range("J131").Value = 792.2023
$ws.Range("K131").Value = 905.45454
$ws.Range("L131").Value = 2376.6069
$ws.Range("M131").Value = 4134.54546
$ws.Range("N131").Value = -12456.6069

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1266.3334
$ws.Range("I16").Value = 1119.6
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1119.6
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -949.5999999999999
$ws.Range("N16").Value = -2340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2572.25
$ws.Range("J122").Value = 2996
$ws.Range("L122").Value = 8988
$ws.Range("N122").Value = -13888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 45461340
$ws.Range("I126").Value = 58831580
$ws.Range("J126").Value = 2520
$ws.Range("K126").Value = 176494740
$ws.Range("L126").Value = 7560
$ws.Range("M126").Value = -176492270
$ws.Range("N126").Value = -12500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1617.8276
$ws.Range("I136").Value = 1600.68
$ws.Range("J136").Value = 1725
$ws.Range("K136").Value = 4802.04
$ws.Range("L136").Value = 5175
$ws.Range("M136").Value = -2252.04
$ws.Range("N136").Value = -10275
